$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "نواقص الأصناف" (shortages) report lost 3 items from the table:
#   row 10 -> PK-MERZ 100MG 30 F.C. TAB
#   row 12 -> TRIVASTAL RETARD 50MG 30 TAB.
#   row 13 -> ZURCAL 40MG 14 GASTRO RESISTANT TAB
# Deleting the corresponding worksheet rows shifts everything below
# up by one slot each time, so delete starting from the top-most row
# so the row indices of the remaining deletions stay correct.
$ws.Rows.Item(10).Delete()
$ws.Rows.Item(11).Delete()
$ws.Rows.Item(11).Delete()

# The grand-total cell (now on row 14, column N) held a static cached
# number rather than a live formula, so it must be refreshed by hand to
# match the sum of the remaining 7 rows' sale-price column (P7:P13).
$ws.Range("N14").Value = 950.095

# Column A is the "م" (serial number) counter, stored as plain literal
# numbers rather than a formula, so deleting rows shifted the old
# values up along with their row instead of renumbering. Restore the
# expected sequential numbering 1-7 for the remaining item rows.
$ws.Range("A10").Value = 4
$ws.Range("A11").Value = 5
$ws.Range("A12").Value = 6
$ws.Range("A13").Value = 7

# Restore the exact row heights of the remaining item rows (10,12,13)
# as shipped in the regenerated report.
$ws.Rows.Item(10).RowHeight = 24.75
$ws.Rows.Item(12).RowHeight = 25.5
$ws.Rows.Item(13).RowHeight = 24.75
